# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (only) Slide Master, currently
#                            holding the "Integral" colour scheme.
#   ppt/theme/theme2.xml  -> bound to the Notes Master, currently
#                            holding the stock "Office Theme" colour scheme.
#
# The target edit swaps the two themes' colour schemes (the font scheme
# and format scheme are already byte-for-byte identical between the two
# theme parts, so only the 12 clrScheme colours actually change).
#
# PowerPoint's object model exposes the Slide Master's theme colours as a
# 12-slot ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) via Master.Theme.ThemeColorScheme.Item(1..12).RGB, so we push
# the "Office Theme" palette (formerly in theme2.xml) onto that theme.

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$themeColors = $master.Theme.ThemeColorScheme

# Target palette (the stock "Office Theme" colours), in
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink order.
$officePalette = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 0; $i -lt $officePalette.Count; $i++) {
    $hex = $officePalette[$i]
    $r = $hex -band 0xFF0000
    $r = $r -shr 16
    $g = $hex -band 0x00FF00
    $g = $g -shr 8
    $b = $hex -band 0x0000FF

    # PowerPoint's RGBColor.RGB takes a COLORREF-style integer: R | (G<<8) | (B<<16)
    $rgbValue = $r + ($g * 256) + ($b * 65536)

    $themeColors.Item($i + 1).RGB = $rgbValue
}

Write-Host "Theme colours swapped to the Office Theme palette."
